# "all SE constraints added (infeasible)"
# Update max_capacity and reorder_level sheets with the new (infeasible)
# spare-parts capacity / reorder-level constraints, and nudge a couple of
# view/selection properties to match what Excel recorded when the sheets
# were last touched.

$wb = $excel.ActiveWorkbook

# --- max_capacity sheet: every S1..S4 x B1/B2 cell becomes 100 -------------
$maxCap = $wb.Worksheets.Item("max_capacity")
$maxCap.Activate() | Out-Null
$maxCap.Range("B2:C5").Value = 100
$maxCap.Range("D11").Select() | Out-Null

# --- reorder_level sheet: every S1..S4 x B1/B2 cell becomes 10 -------------
$reorder = $wb.Worksheets.Item("reorder_level")
$reorder.Activate() | Out-Null
$reorder.Range("B2:C5").Value = 10
$reorder.Range("C2").Select() | Out-Null

# --- spare_parts sheet: only the scrolled view position changed ------------
$spareParts = $wb.Worksheets.Item("spare_parts")
$spareParts.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 11
$win.ScrollRow = 1
$spareParts.Range("G9").Select() | Out-Null

# --- restore the originally active sheet/tab -------------------------------
$maxCap.Activate() | Out-Null
$maxCap.Range("D11").Select() | Out-Null

# Scroll the sheet tab strip so one sheet earlier is the first visible tab,
# matching the workbookView firstSheet shift recorded in the saved file.
$win = $excel.ActiveWindow
$win.ScrollWorkbookTabs(1, -1) | Out-Null
